$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2022-12-23) is inserted at row 318, pushing the
# existing rows 318-356 down to 319-357 (dimension grows from R356 to R357).
$ws.Rows("318:318").Insert()

$fecha = Get-Date -Year 2022 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(318, 1).Value = 9
$ws.Cells.Item(318, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(318, 3).Value = "Metropolitana"
$ws.Cells.Item(318, 4).Value = $fecha
$ws.Cells.Item(318, 5).Value = 13
$ws.Cells.Item(318, 6).Value = 300000001
$ws.Cells.Item(318, 7).Value = "Rabanito"
$ws.Cells.Item(318, 8).Value = "Sin especificar"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 7000
$ws.Cells.Item(318, 11).Value = 3000
$ws.Cells.Item(318, 12).Value = 3000
$ws.Cells.Item(318, 13).Value = 3000
$ws.Cells.Item(318, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(318, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(318, 16).Value = 30
$ws.Cells.Item(318, 17).Value = 100
$ws.Cells.Item(318, 18).Value = "Hortaliza"
